{"js": "// Update the worksheet date and the 25 division problems/answers in the table.\n// The edits are applied positionally (by paragraph index in document order),\n// because a couple of the \"before\" strings are not unique in the document\n// (e.g. \"400\u00f72=200, 0\" appears twice but maps to two different results).\nconst edits = [\n  { index: 0, before: \"2026-02-24 Tuesday\", after: \"2026-02-25 Wednesday\" },\n  { index: 1, before: \"925\u00f77=132, 1\", after: \"543\u00f74=135, 3\" },\n  { index: 2, before: \"658\u00f79=73, 1\", after: \"475\u00f74=118, 3\" },\n  { index: 3, before: \"706\u00f79=78, 4\", after: \"765\u00f72=382, 1\" },\n  { index: 4, before: \"360\u00f79=40, 0\", after: \"558\u00f73=186, 0\" },\n  { index: 5, before: \"144\u00f72=72, 0\", after: \"567\u00f77=81, 0\" },\n  { index: 21, before: \"880\u00f79=97, 7\", after: \"854\u00f75=170, 4\" },\n  { index: 22, before: \"344\u00f79=38, 2\", after: \"816\u00f72=408, 0\" },\n  { index: 23, before: \"478\u00f76=79, 4\", after: \"834\u00f73=278, 0\" },\n  { index: 24, before: \"826\u00f78=103, 2\", after: \"262\u00f72=131, 0\" },\n  { index: 25, before: \"975\u00f79=108, 3\", after: \"753\u00f76=125, 3\" },\n  { index: 41, before: \"662\u00f76=110, 2\", after: \"635\u00f77=90, 5\" },\n  { index: 42, before: \"276\u00f72=138, 0\", after: \"465\u00f77=66, 3\" },\n  { index: 43, before: \"112\u00f73=37, 1\", after: \"882\u00f75=176, 2\" },\n  { index: 44, before: \"682\u00f74=170, 2\", after: \"457\u00f73=152, 1\" },\n  { index: 45, before: \"400\u00f72=200, 0\", after: \"108\u00f73=36, 0\" },\n  { index: 61, before: \"217\u00f73=72, 1\", after: \"506\u00f73=168, 2\" },\n  { index: 62, before: \"720\u00f72=360, 0\", after: \"925\u00f73=308, 1\" },\n  { index: 63, before: \"471\u00f75=94, 1\", after: \"273\u00f77=39, 0\" },\n  { index: 64, before: \"413\u00f78=51, 5\", after: \"647\u00f75=129, 2\" },\n  { index: 65, before: \"785\u00f74=196, 1\", after: \"158\u00f73=52, 2\" },\n  { index: 81, before: \"741\u00f75=148, 1\", after: \"295\u00f74=73, 3\" },\n  { index: 82, before: \"119\u00f75=23, 4\", after: \"347\u00f73=115, 2\" },\n  { index: 83, before: \"641\u00f79=71, 2\", after: \"147\u00f79=16, 3\" },\n  { index: 84, before: \"400\u00f72=200, 0\", after: \"240\u00f74=60, 0\" },\n  { index: 85, before: \"994\u00f75=198, 4\", after: \"276\u00f72=138, 0\" },\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst maxIndex = Math.max(...edits.map((e) => e.index));\nif (paragraphs.items.length <= maxIndex) {\n  throw new Error(\n    `Expected at least ${maxIndex + 1} paragraphs, found ${paragraphs.items.length}`\n  );\n}\n\nfor (const edit of edits) {\n  const paragraph = paragraphs.items[edit.index];\n  const current = paragraph.text;\n  if (current !== edit.before) {\n    throw new Error(\n      `Paragraph ${edit.index}: expected \"${edit.before}\" but found \"${current}\"`\n    );\n  }\n  paragraph.insertText(edit.after, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date (first paragraph) and the 25 division\n# problems/answers that live in the single 20-row x 5-column table.\n# Addressing is done by paragraph index / table cell coordinates rather\n# than by text search-and-replace, because a couple of the \"before\"\n# strings (e.g. \"400\u00f72=200, 0\") are not unique in the document and map\n# to two different replacement values depending on position.\n\n$d = $word.ActiveDocument\n\n# --- Title / date line -------------------------------------------------\n$dateParagraph = $d.Paragraphs.Item(1)\nif ($dateParagraph.Range.Text.TrimEnd([char]13, [char]7) -ne \"2026-02-24 Tuesday\") {\n    throw \"Unexpected text in date paragraph: $($dateParagraph.Range.Text)\"\n}\n$dateParagraph.Range.Text = \"2026-02-25 Wednesday\"\n\n# --- Table of division problems ----------------------------------------\n$table = $d.Tables.Item(1)\n\n$edits = @(\n    @{ Row = 1;  Col = 1; Before = \"925\u00f77=132, 1\";  After = \"543\u00f74=135, 3\" },\n    @{ Row = 1;  Col = 2; Before = \"658\u00f79=73, 1\";   After = \"475\u00f74=118, 3\" },\n    @{ Row = 1;  Col = 3; Before = \"706\u00f79=78, 4\";   After = \"765\u00f72=382, 1\" },\n    @{ Row = 1;  Col = 4; Before = \"360\u00f79=40, 0\";   After = \"558\u00f73=186, 0\" },\n    @{ Row = 1;  Col = 5; Before = \"144\u00f72=72, 0\";   After = \"567\u00f77=81, 0\" },\n\n    @{ Row = 5;  Col = 1; Before = \"880\u00f79=97, 7\";   After = \"854\u00f75=170, 4\" },\n    @{ Row = 5;  Col = 2; Before = \"344\u00f79=38, 2\";   After = \"816\u00f72=408, 0\" },\n    @{ Row = 5;  Col = 3; Before = \"478\u00f76=79, 4\";   After = \"834\u00f73=278, 0\" },\n    @{ Row = 5;  Col = 4; Before = \"826\u00f78=103, 2\";  After = \"262\u00f72=131, 0\" },\n    @{ Row = 5;  Col = 5; Before = \"975\u00f79=108, 3\";  After = \"753\u00f76=125, 3\" },\n\n    @{ Row = 9;  Col = 1; Before = \"662\u00f76=110, 2\";  After = \"635\u00f77=90, 5\" },\n    @{ Row = 9;  Col = 2; Before = \"276\u00f72=138, 0\";  After = \"465\u00f77=66, 3\" },\n    @{ Row = 9;  Col = 3; Before = \"112\u00f73=37, 1\";   After = \"882\u00f75=176, 2\" },\n    @{ Row = 9;  Col = 4; Before = \"682\u00f74=170, 2\";  After = \"457\u00f73=152, 1\" },\n    @{ Row = 9;  Col = 5; Before = \"400\u00f72=200, 0\";  After = \"108\u00f73=36, 0\" },\n\n    @{ Row = 13; Col = 1; Before = \"217\u00f73=72, 1\";   After = \"506\u00f73=168, 2\" },\n    @{ Row = 13; Col = 2; Before = \"720\u00f72=360, 0\";  After = \"925\u00f73=308, 1\" },\n    @{ Row = 13; Col = 3; Before = \"471\u00f75=94, 1\";   After = \"273\u00f77=39, 0\" },\n    @{ Row = 13; Col = 4; Before = \"413\u00f78=51, 5\";   After = \"647\u00f75=129, 2\" },\n    @{ Row = 13; Col = 5; Before = \"785\u00f74=196, 1\";  After = \"158\u00f73=52, 2\" },\n\n    @{ Row = 17; Col = 1; Before = \"741\u00f75=148, 1\";  After = \"295\u00f74=73, 3\" },\n    @{ Row = 17; Col = 2; Before = \"119\u00f75=23, 4\";   After = \"347\u00f73=115, 2\" },\n    @{ Row = 17; Col = 3; Before = \"641\u00f79=71, 2\";   After = \"147\u00f79=16, 3\" },\n    @{ Row = 17; Col = 4; Before = \"400\u00f72=200, 0\";  After = \"240\u00f74=60, 0\" },\n    @{ Row = 17; Col = 5; Before = \"994\u00f75=198, 4\";  After = \"276\u00f72=138, 0\" }\n)\n\nforeach ($edit in $edits) {\n    $cell = $table.Cell($edit.Row, $edit.Col)\n    $cellRange = $cell.Range\n    $currentText = $cellRange.Text.TrimEnd([char]13, [char]7)\n    if ($currentText -ne $edit.Before) {\n        throw \"Cell ($($edit.Row),$($edit.Col)): expected '$($edit.Before)' but found '$currentText'\"\n    }\n    $cellRange.Text = $edit.After\n}\n"}
